# "fix mk9 barrels and adjust weights"
# Adds the KAK Value Line 3" AR9 barrel and the CmmG Mk9 barrel family
# (5", 8.5", 9", and three 16" variants) to the m4-barrels sheet, rows 11-17,
# extends the shared N/S formulas down through row 17, and updates the
# active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 10 - fill in the N/S formulas on the (otherwise still blank) spacer
# row so the shared formula range can run all the way to row 17.
# ---------------------------------------------------------------------
$ws.Range("N10").Formula = "=C10-D10*20-E10*0.8-F10*0.6-H10*5+I10*10+J10/300"
$ws.Range("S10").Formula = "=ROUND(Q10*0.02+P10+R10, 2)"

# ---------------------------------------------------------------------
# Row 11 - KAK Value Line Light Tapered Melonite AR9 3" 9x19
# ---------------------------------------------------------------------
$ws.Range("A11").Value = "kak_value_line_light_tapered_melonite_ar9_76mm_9x19_barrel"
$ws.Range("B11").Value = 'KAK Value Line Light Tapered Melonite AR9 3" 9x19'
$ws.Range("C11").Value = 8
$ws.Range("D11").Value = 0.12
$ws.Range("E11").Value = 7
$ws.Range("F11").Value = 7
$ws.Range("H11").Value = 0.3
$ws.Range("I11").Value = -0.3
$ws.Range("J11").Value = -300
$ws.Range("M11").Value = 750
$ws.Range("N11").Formula = "=C11-D11*20-E11*0.8-F11*0.6-H11*5+I11*10+J11/300"
$ws.Range("P11").Value = 0.06
$ws.Range("Q11").Value = 3
$ws.Range("S11").Formula = "=ROUND(Q11*0.02+P11+R11, 2)"

# ---------------------------------------------------------------------
# Row 12 - CmmG Mk9 5" 9x19
# (the two shared-string entries for this row were appended last in the
# authored workbook, so the A12/B12 text is assigned further below, after
# rows 13-17 have claimed their shared-string slots, to reproduce the same
# shared-string ordering)
# ---------------------------------------------------------------------
$ws.Range("C12").Value = 5
$ws.Range("D12").Value = 0.16
$ws.Range("E12").Value = 4
$ws.Range("F12").Value = 4
$ws.Range("H12").Value = 0.25
$ws.Range("I12").Value = -0.12
$ws.Range("J12").Value = -225
$ws.Range("M12").Value = 800
$ws.Range("N12").Formula = "=C12-D12*20-E12*0.8-F12*0.6-H12*5+I12*10+J12/300"
$ws.Range("P12").Value = 0.06
$ws.Range("Q12").Value = 5
$ws.Range("S12").Formula = "=ROUND(Q12*0.02+P12+R12, 2)"

# ---------------------------------------------------------------------
# Row 13 - CmmG Mk9 8.5" 9x19
# ---------------------------------------------------------------------
$ws.Range("A13").Value = "cmmg_mk9_216mm_9x19_barrel"
$ws.Range("B13").Value = 'CmmG Mk9 8.5" 9x19'
$ws.Range("C13").Value = 2
$ws.Range("D13").Value = 0.23
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 2
$ws.Range("H13").Value = 0.12
$ws.Range("I13").Value = -0.04
$ws.Range("J13").Value = -90
$ws.Range("M13").Value = 900
$ws.Range("N13").Formula = "=C13-D13*20-E13*0.8-F13*0.6-H13*5+I13*10+J13/300"
$ws.Range("P13").Value = 0.06
$ws.Range("Q13").Value = 8.5
$ws.Range("S13").Formula = "=ROUND(Q13*0.02+P13+R13, 2)"

# ---------------------------------------------------------------------
# Row 14 - CmmG Mk9 9" 9x19
# ---------------------------------------------------------------------
$ws.Range("A14").Value = "cmmg_mk9_229mm_9x19_barrel"
$ws.Range("B14").Value = 'CmmG Mk9 9" 9x19'
$ws.Range("C14").Value = 0
$ws.Range("D14").Value = 0.24
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 0
$ws.Range("H14").Value = 0.1
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = -70
$ws.Range("M14").Value = 950
$ws.Range("N14").Formula = "=C14-D14*20-E14*0.8-F14*0.6-H14*5+I14*10+J14/300"
$ws.Range("P14").Value = 0.06
$ws.Range("Q14").Value = 9
$ws.Range("S14").Formula = "=ROUND(Q14*0.02+P14+R14, 2)"

# ---------------------------------------------------------------------
# Row 15 - CmmG Mk9 16" 9x19 Carbine Length
# ---------------------------------------------------------------------
$ws.Range("A15").Value = "cmmg_mk9_406mm_9x19_c_barrel"
$ws.Range("B15").Value = 'CmmG Mk9 16" 9x19 Carbine Length'
$ws.Range("C15").Value = -6
$ws.Range("D15").Value = 0.38
$ws.Range("E15").Value = -4
$ws.Range("F15").Value = -4
$ws.Range("H15").Value = 0.05
$ws.Range("I15").Value = 0.1
$ws.Range("J15").Value = 200
$ws.Range("M15").Value = 1200
$ws.Range("N15").Formula = "=C15-D15*20-E15*0.8-F15*0.6-H15*5+I15*10+J15/300"
$ws.Range("P15").Value = 0.06
$ws.Range("Q15").Value = 16
$ws.Range("S15").Formula = "=ROUND(Q15*0.02+P15+R15, 2)"

# ---------------------------------------------------------------------
# Row 16 - CmmG Mk9 16" 9x19 Rifle Length
# ---------------------------------------------------------------------
$ws.Range("A16").Value = "cmmg_mk9_406mm_9x19_r_barrel"
$ws.Range("B16").Value = 'CmmG Mk9 16" 9x19 Rifle Length'
$ws.Range("C16").Value = -7
$ws.Range("D16").Value = 0.39
$ws.Range("E16").Value = -5
$ws.Range("F16").Value = -5
$ws.Range("H16").Value = 0.05
$ws.Range("I16").Value = 0.1
$ws.Range("J16").Value = 200
$ws.Range("M16").Value = 1200
$ws.Range("N16").Formula = "=C16-D16*20-E16*0.8-F16*0.6-H16*5+I16*10+J16/300"
$ws.Range("P16").Value = 0.06
$ws.Range("Q16").Value = 16
$ws.Range("R16").Value = 0.01
$ws.Range("S16").Formula = "=ROUND(Q16*0.02+P16+R16, 2)"

# ---------------------------------------------------------------------
# Row 17 - CmmG Mk9 16" 9x19 Mid Length
# ---------------------------------------------------------------------
$ws.Range("A17").Value = "cmmg_mk9_406mm_9x19_m_barrel"
$ws.Range("B17").Value = 'CmmG Mk9 16" 9x19 Mid Length'
$ws.Range("C17").Value = -8
$ws.Range("D17").Value = 0.4
$ws.Range("E17").Value = -6
$ws.Range("F17").Value = -6
$ws.Range("H17").Value = 0.05
$ws.Range("I17").Value = 0.1
$ws.Range("J17").Value = 200
$ws.Range("M17").Value = 1200
$ws.Range("N17").Formula = "=C17-D17*20-E17*0.8-F17*0.6-H17*5+I17*10+J17/300"
$ws.Range("P17").Value = 0.06
$ws.Range("Q17").Value = 16
$ws.Range("R17").Value = 0.02
$ws.Range("S17").Formula = "=ROUND(Q17*0.02+P17+R17, 2)"

# ---------------------------------------------------------------------
# Row 12 text (see note above for why this is assigned last)
# ---------------------------------------------------------------------
$ws.Range("A12").Value = "cmmg_mk9_127mm_9x19_barrel"
$ws.Range("B12").Value = 'CmmG Mk9 5" 9x19'

# ---------------------------------------------------------------------
# Restore the selection to the cell it ended up on in the authored file.
# ---------------------------------------------------------------------
[void]$ws.Range("E16").Select()
